$wb = $excel.ActiveWorkbook

# --- Workbook window position/size (cosmetic) ---
$excel.Windows.Item(1).Left = 540
$excel.Windows.Item(1).Top = 660
$excel.Windows.Item(1).Width = 12450
$excel.Windows.Item(1).Height = 15165

# --- Insert new "adj_perc" column into "Contracting by Year" ---
$ws = $wb.Worksheets.Item("Contracting by Year")

# Shift B:G -> C:H and create a fresh column B
$ws.Columns("B:B").Insert()

$ws.Range("B1").Value = "adj_perc"

$ws.Range("B2").Value = 49
$ws.Range("B3").Value = 48
$ws.Range("B4").Value = 56
$ws.Range("B5").Value = 58
$ws.Range("B6").Value = 59
$ws.Range("B7").Value = 63
$ws.Range("B8").Value = 61
$ws.Range("B9").Value = 60
$ws.Range("B10").Value = 53

# Match the new column's look to the rest of the plain data columns
$ws.Range("B2:B10").Font.Name = "Arial"
$ws.Range("B2:B10").Font.Size = 10
$ws.Range("B2:B10").NumberFormat = "General"
